$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; this shifts the existing rows 54-122
# down to 55-123 (all their data and styles move down with them).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44650
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100109
$ws.Range("H54").Value = "Uva"
$ws.Range("I54").Value = 100109001
$ws.Range("J54").Value = "Uva"
$ws.Range("K54").Value = "Thompson seedless"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 130
$ws.Range("N54").Value = 10000
$ws.Range("O54").Value = 11000
$ws.Range("P54").Value = 10385
$ws.Range("Q54").Value = "$/bandeja 18 kilos"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 577
$ws.Range("T54").Value = 18
